$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add formula to M26 (Capital row), computing 1.5x the L26 value
$ws.Range("M26").Formula = "=L26*1.5"

# Update the sheet view: enable right-to-left, scroll to H10, and select M27
$win = $excel.ActiveWindow
$win.DisplayRightToLeft = $true
$win.ScrollRow = 10
$win.ScrollColumn = 8
$ws.Range("M27").Select() | Out-Null
